# Apply cryptos-list update (price + volume refresh, RenderToken/Kaspa row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.445.47"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.95%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.212.00"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.16%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.65"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.616"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.97"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.06%  "

$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.600"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.29%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.32"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.91%  "

$ws.Range("E11").Value = "  -2.39%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.96"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.72%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.85"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.86%  "

$ws.Range("E14").Value = "  -2.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.544.19"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.14%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.68"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.38%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.213.97"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.85%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.800"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.52%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.336.20"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.91%  "

$ws.Range("E20").Value = "  -0.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.77"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.91"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.75%  "

$ws.Range("E23").Value = "  -9.98%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "228.60"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.10"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.02%  "

$ws.Range("E26").Value = "  -0.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.92"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.37"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -6.92%  "

$ws.Range("E29").Value = "  -2.20%  "

$ws.Range("E30").Value = "  -0.93%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.69"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.92%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.16"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.79"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +12.54%  "

$ws.Range("E34").Value = "  -0.90%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.37"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.83%  "

$ws.Range("E36").Value = "  -2.27%  "

$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.108"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.76%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.40"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0322"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +6.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.61"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.63%  "

$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("E42").Value = "  -3.65%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "60.46"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.196"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.54"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.78%  "

$ws.Range("E46").Value = "  -2.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "99.37"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.06%  "

$ws.Range("E48").Value = "  -2.93%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.30"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.30%  "

$ws.Range("E50").Value = "  -2.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.420"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +13.75%  "
